# Update the "Förändrad" (Changed) date column (C) for rows 2-24
# from serial date 45209 (2023-10-10) to 45210 (2023-10-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 24; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45209) {
        $cell.Value2 = 45210
    }
}
